$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update "Valor Mora" total and "Cant. Periodos" count
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 2553600
$ws.Range("F13").Value = 19

# ---------------------------------------------------------------------------
# 2) Re-sort the "Periodo Mora" column (E16:E33) from descending to ascending
#    order (2402 .. 2507).
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "2402"
$ws.Range("E17").Value = "2403"
$ws.Range("E18").Value = "2404"
$ws.Range("E19").Value = "2405"
$ws.Range("E20").Value = "2406"
$ws.Range("E21").Value = "2407"
$ws.Range("E22").Value = "2408"
$ws.Range("E23").Value = "2409"
$ws.Range("E24").Value = "2410"
$ws.Range("E25").Value = "2411"
$ws.Range("E26").Value = "2412"
$ws.Range("E27").Value = "2501"
$ws.Range("E28").Value = "2502"
$ws.Range("E29").Value = "2503"
$ws.Range("E30").Value = "2504"
$ws.Range("E31").Value = "2505"
$ws.Range("E32").Value = "2506"
$ws.Range("E33").Value = "2507"

# ---------------------------------------------------------------------------
# 3) Insert a new data row (34) for period 2508. Before inserting, grab the
#    distinctive "bottom border / last row" formatting that row 33 currently
#    has, so the new row keeps that look; then restyle row 33 back to the
#    regular "middle of table" formatting (copied from row 32).
# ---------------------------------------------------------------------------
$ws.Range("B33:J33").Copy() | Out-Null
$ws.Range("B34:J34").EntireRow.Insert()
$ws.Range("B34:J34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B32:J32").Copy() | Out-Null
$ws.Range("B33:J33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "11511691"
$ws.Range("D34").Value = "JAIRO ALFREDO JIMENEZ BARON"
$ws.Range("E34").Value = "2508"
$ws.Range("F34").Value = 134400
$ws.Range("G34").Value = 3360000

# Note: the single row-insert above (at row 34) already shifts everything
# below it down by one, so the old "signature line" row (38) now lives at
# row 39 and the old "name / signature" labels row (39) now lives at row 40
# - exactly the target layout, with no further row manipulation required.

Write-Host "Edit complete"
